# EXPORT BOT6 Laufzettel - add a second sheet "Tabelle1" that duplicates the
# content of the original "Beispiel 1 DIN A4 Seite" sheet (fileReader fix /
# aktuelle Import Dateien), and switch the active tab to it.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Duplicate the existing sheet and put the copy at the end ---------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2.Name = "Tabelle1"

# --- 2. Drop the left-over helper cells (G13:J15) that only existed to   --
#        pad sheet1's dimension - the new sheet's used range is A1:F35.
$ws2.Range("G13:J15").Clear()

# --- 3. Re-flow row heights / fonts like the freshly generated sheet -----
#        (smaller default row height, long-text rows auto growing).
$ws2.Range("A1:F35").EntireRow.RowHeight = 15.75
$ws2.Range("A8:F8").EntireRow.RowHeight = 78.75
$ws2.Range("A26:F26").EntireRow.RowHeight = 78.75

# --- 4. Add the thin blank separator rows between the 4 student blocks ---
$blankRows = @(9, 18, 27)
foreach ($r in $blankRows) {
    $rowRange = $ws2.Range("A" + $r + ":F" + $r)
    $rowRange.Font.Size = 12
    $rowRange.Font.Bold = $false
    $rowRange.VerticalAlignment = -4108
    $ws2.Range("C" + $r).HorizontalAlignment = -4108
    $ws2.Range("F" + $r).HorizontalAlignment = -4108
}

# --- 5. Narrower, auto-fit-like column widths on the new sheet -----------
$ws2.Columns.Item(1).ColumnWidth = 18.166666666666668
$ws2.Columns.Item(2).ColumnWidth = 11.166666666666666
$ws2.Columns.Item(3).ColumnWidth = 12.451822916666666

# --- 6. Page margins matching the regenerated sheet -----------------------
$ws2.PageSetup.LeftMargin = 50.4
$ws2.PageSetup.RightMargin = 50.4
$ws2.PageSetup.TopMargin = 56.6929133999999
$ws2.PageSetup.BottomMargin = 56.6929133999999
$ws2.PageSetup.HeaderMargin = 21.6
$ws2.PageSetup.FooterMargin = 21.6

# --- 7. Selections: old sheet keeps cursor at K23, new sheet at K17 ------
$ws1.Range("K23").Select()
$ws2.Range("K17").Select()

# --- 8. New sheet becomes the active / displayed tab ----------------------
$ws2.Select()
